$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.991.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.74%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.333.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.92%  "

$ws.Range("E4").Value = "  -0.50%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.91%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.63%  "

$ws.Range("E7").Value = "  +3.87%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.51%  "

$ws.Range("E9").Value = "  +8.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.07"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0816"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.48"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.29%  "

$ws.Range("E13").Value = "  +0.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.688.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.335.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.62%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.839"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.00%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "46.857.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +19.64%  "

$ws.Range("E20").Value = "  +4.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.74%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "252.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.25%  "

$ws.Range("E24").Value = "  +3.76%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.29%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "42.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +16.29%  "

$ws.Range("E28").Value = "  +1.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.94%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.58%  "

$ws.Range("E32").Value = "  +7.71%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "146.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.15%  "

$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.116"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.21%  "

$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.26%  "

$ws.Range("E37").Value = "  +3.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.25%  "

$ws.Range("E39").Value = "  +9.39%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0311"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.86%  "

$ws.Range("E41").Value = "  +4.92%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.72%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +16.53%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.805.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +15.47%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "74.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +11.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.195"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "99.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.62%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "55.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.68%  "

$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.24%  "
